# KP-11725 D: Extension of questionnaire's translation files
# Insert a new "Variable" column (with constant value "c1" for every data
# row) right after the "Entity Id" column on both sheets, and update the
# active sheet / selection to reflect where the edit left the cursor.

$wb = $excel.ActiveWorkbook

# --- "Translations" sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Translations")
$ws1.Columns("B:B").Insert()

$ws1.Range("B1").Value = "Variable"
$ws1.Range("B2").Value = "c1"
$ws1.Range("B3").Value = "c1"
$ws1.Range("B4").Value = "c1"
$ws1.Range("B5").Value = "c1"

# --- "@@_question" sheet ----------------------------------------------------
$ws2 = $wb.Worksheets.Item("@@_question")
$ws2.Columns("B:B").Insert()

$ws2.Range("B1").Value = "Variable"
$ws2.Range("B2").Value = "c1"

# --- Selection / active sheet, matching where the edit left the cursor -----
$ws2.Activate()
$ws2.Range("B2").Select() | Out-Null

$ws1.Activate()
$ws1.Range("B6").Select() | Out-Null
